$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset column R ("backup") back to 0 for the rows that were recomputed.
$rowsToZero = @(55,63,68,70,81,87,113,114,126,133,137,150,156,162,168,181,200,203,211,218,231,238)
foreach ($r in $rowsToZero) {
    $ws.Cells.Item($r, 18).Value = 0
}

# Row 242's "isPivot" flag (column O) is now set.
$ws.Cells.Item(242, 15).Value = 1

# Append the two new weekly bars that were broken out of stock.yaml.
$ws.Cells.Item(246, 1).Value = 45460
$ws.Cells.Item(246, 2).Value = 1018.200012207031
$ws.Cells.Item(246, 3).Value = 1046
$ws.Cells.Item(246, 4).Value = 1002.900024414062
$ws.Cells.Item(246, 5).Value = 1012.299987792969
$ws.Cells.Item(246, 6).Value = 1012.299987792969
$ws.Cells.Item(246, 7).Value = 12396390
$ws.Cells.Item(246, 8).Value = 2024
$ws.Cells.Item(246, 9).Value = 6
$ws.Cells.Item(246, 10).Value = 17
$ws.Cells.Item(246, 11).Value = 0
$ws.Cells.Item(246, 12).Value = 0
$ws.Cells.Item(246, 13).Value = 0
$ws.Cells.Item(246, 14).Value = 25
$ws.Cells.Item(246, 15).Value = 0
$ws.Cells.Item(246, 16).Value = 0
$ws.Cells.Item(246, 17).Value = 0
$ws.Cells.Item(246, 18).Value = "'"
$ws.Cells.Item(246, 18).Style = "Normal"

$ws.Cells.Item(247, 1).Value = 45467
$ws.Cells.Item(247, 2).Value = 1020
$ws.Cells.Item(247, 3).Value = 1037.949951171875
$ws.Cells.Item(247, 4).Value = 983.3499755859375
$ws.Cells.Item(247, 5).Value = 989.25
$ws.Cells.Item(247, 6).Value = 989.25
$ws.Cells.Item(247, 7).Value = 16612019
$ws.Cells.Item(247, 8).Value = 2024
$ws.Cells.Item(247, 9).Value = 6
$ws.Cells.Item(247, 10).Value = 24
$ws.Cells.Item(247, 11).Value = 0
$ws.Cells.Item(247, 12).Value = 0
$ws.Cells.Item(247, 13).Value = 0
$ws.Cells.Item(247, 14).Value = 26
$ws.Cells.Item(247, 15).Value = 0
$ws.Cells.Item(247, 16).Value = 0
$ws.Cells.Item(247, 17).Value = 0
$ws.Cells.Item(247, 18).Value = "'"
$ws.Cells.Item(247, 18).Style = "Normal"

# Match the date-time number format used by the rest of column A.
$ws.Range("A246:A247").NumberFormat = "YYYY-MM-DD HH:MM:SS"
